$wb = $excel.ActiveWorkbook

# --- Update shared "per-play" stat strings on the YDS and ST sheets (Week 16 logged) ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsST = $wb.Worksheets.Item("ST")

$s4 = "5 3 5 3 3 5 2 5 1 7 0 6 6 1 8 2 7 3 2 10 3 8 10 4 2 8 4 6 0 3 1 3 1 5 -1 4 3 0 2 3 1 4 9 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 1 17 4 8 4 4 4 6 4 -1 13 13 3 6 1 4 5 5 4 6 2 4 0 9 4 3 6 0 2 3 2 2 20 5 6 2 7 0 1 9 2 2 0 1 -1 3 15 2 4 5 12 3 9 4 11 6 8 12 2 -1 1 12 19 12 9 5 3 1 2 -3 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 5 6 0 3 5 8 2 1 3 4 8 6 7 10 7 6 3 2 5 10 3 3 4 6 2 7 -1 3 4 9 2 2 1 3 -3 2 -1 4 7 3 7 -3 4 2 2 1 11 6 12 4 1 4 5 8 8 1 2 3 1 11 8 3 28 3 4 2 2 6 6 5 4 3 1 3 4 9 -4 1 9 0 3 47 4 2 12 2 3 1 5 -1 8 2 1 1 2 9 1 -2 2 0 7 21 1 4 3 1 3 6 2 3 5 2 -1 6 6 3 8 0 3 12 30 2 3 2 6 0 17 5 3 4 5 55 0 4 4 4 -2 9 1 2 3 -1 3 1 1 2 0 11 14 7 2 -3 16 -1 12 -1 0 -1"
$s5 = "10 19 28 16 8 5 27 5 19 17 2 10 47 20 2 2 11 10 1 5 7 9 7 9 7 3 11 8 4 6 20 24 20 25 20 16 16 9 5 8 1 9 5 26 17 21 3 8 1 2 9 10 26 5 12 6 8 5 1 15 17 4 3 7 7 0 17 5 26 3 20 20 32 14 8 20 15 11 6 7 17 11 -4 23 13 31 17 8 -2 11 8 5 5 21 -1 9 9 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 9 26 13 10 10 18 6 15 62 6 10 10 11 12 3 4 11 5 12 6 14 12 16 4 15 7 34 19 9 22 9 22 16 10 2 8 8 11 3 9 23 -4 10 12 14 18 10 3 6 -4 5 4 4 9 15 0 8 5 4 8 10 2 27 10 -4 6 4 21 4 11 9 10 46 2 9 14 0 9 8 9 7 12 0 34 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 2 3 29 9 8 3 6 7 11 4 5 14 15 8 6 2 15 4 0 4 3 40 10 8 15 3 13 7 10 8 6 6 6 7 4 35 7 12 16 3 13 9 17 14 14 18 4 5 8 3 21 5 3 20 9 0 0 4 -1 10 8 7 7 15 9 26 5 4 11 10 32 9 2 9 3 8 6 13 15 5 11 0 6 3 15 7 4 3 25 36 3 0 15 6 17 12 27 1 9 1 7 11 5 7 18 1 5 11 14 6 6 6 36 5 3 6 15 5 7 5 9 20 8 16 23 18 5 13 16 13 7 5 20 5 6 24 4 18 0 8 5 11 10 6 -7 14 58 5 7 7 -2 3 23 5 8 1 5 5 6 14 19 7 5 8 5 16 6 10 3 6 18 11 13 4 3 4 13 7 19 7 1 10 22 62 13 23 11 12 7 10"
$s6 = "4 2 4 0 2 7 1 0 3 -1 13 3 0 2 4 1 4 11 7 -1 0 2 3 10 0 0 8 2 0 8 1 -3 1 8 6 0 3 6 4 0 4 4 2 5 1 4 -2 0 2 11 2 3 1 4 1 5 15 0 6 -4 0 0 0 -1 1 2 4 0 3 1 13 7 4 5 -2 11 1 1 -1 5 6 2 5 -2 9 11 23 14 7 3 -1 2 0 29 1 1 6 4 0 6 13 5 12 6 0 6 0 8 8 1 2 -1 4 -1 -1 1 1 8 2 6 16 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 1 3 1 5 -1 -8 10 5 2 3 15 4 4 -5 -1 -3 12 9 1 1 1 2 3 4 7 3 2 4 2 -1 5 0 3 2 8 6 0 1 3 16 1 2 12 0 9 6 0 4 -2 0 5 12 5 1 16 3 2 5 15 5 10 15 1 3 4 7 39 17 6 5 1 5 11 0 4 2 0 4 0 3 1 0 2 3 0 14 0 14 -3 10 22 1 29 18 4 11 5 18 23 2 0 7 8 3 8 -1 4 -3 5 -2 1 1 9 2 2 6 6 3 -1 3 8 0 0 7 3 -1 4 1 1 -1 3 -3 1 0 -5 33 2 -2 1 6 2 2 -2 0 4 5 0 3 3 4 11"
$s7 = "28 5 16 11 8 12 7 8 7 22 2 5 22 6 8 3 4 13 -5 4 11 7 10 5 16 6 14 6 4 3 7 10 9 6 21 2 3 6 11 9 31 10 0 9 5 9 4 5 23 24 4 8 24 7 6 1 14 8 3 8 14 4 23 12 8 7 22 5 7 8 1 2 4 10 2 6 7 6 12 22 11 5 6 6 6 9 2 11 21 22 2 75 0 20 5 40 10 5 17 6 4 10 3 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 21 4 34 9 12 10 5 24 5 3 12 5 2 21 6 10 8 4 23 10 2 7 1 0 11 5 21 23 5 6 4 16 2 1 6 -3 24 25 6 6 5 5 5 3 6 3 17 4 12 4 10 22 20 6 22 12 -6 11 8 5 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 8 14 13 4 8 11 1 20 22 7 9 8 14 1 15 6 18 21 9 6 16 6 5 6 6 9 12 9 4 3 5 5 0 11 21 4 28 8 11 8 -1 7 2 4 5 5 8 4 9 8 62 10 14 3 3 15 9 5 12 22 9 23 -3 15 5 15 11 11 7 8 7 10 0 6 6 12 19 29 7 10 9 8 9 4 -3 14 2 19 6 9 4 16 6 20 6 17 8 8 8 2 7 30 1 10 11 5 3 25 24 21 9 5 6 10 12 6 1 13 9 -2 10 2 4 10 15 9 12 12 3 4 4 4 2 15 1 24 3 8 9 7 40 6 7 33 8 4 3 17 5 7 14 8 10 5 15 63 0 4 6 2 9 19 8 6 13 2 18 14 6 6 17 6"
$s46 = "51 50 44 56 76 64 67 74 73 64 64 69 67 69 66 65 60"
$s47 = "13 19 14 45 27 0 25 33 36 28 19 72 30 29 27 25 26"
$s49 = "65 42 49 41 48 41 43 57 38 53 39 15 55 38 34 43 41 41 44 42 49 42 46 54 32 40 43 46 35 41 55 40 14 54 55 60 41 32 55 20 39 35 47 40 53 48 45 39 42 30 39 45 40 38"
$s50 = "0 3 12 0 6 0 5 0 9 0 0 0 0 0 11 6 0 0 0 8 0 13 0 0 0 12 10 0 0 0 0 0 0 9 13 7 0 0 4 0 0 0 5 0 34 0 0 4 8 0 0 16 6 0"
$s51 = "0 0 10 0 0 0 0 0 6 11 3 2 9 10 0 0 0 4 9 20 6 5"

$wsYDS.Range("B2").Value = $s4
$wsYDS.Range("B3").Value = $s5
$wsYDS.Range("C2").Value = $s6
$wsYDS.Range("C3").Value = $s7
$wsST.Range("B4").Value = $s46
$wsST.Range("B5").Value = $s47
$wsST.Range("D3").Value = $s49
$wsST.Range("D4").Value = $s50
$wsST.Range("D5").Value = $s51

# --- Update aggregate numeric totals across OFF, DEF, ST, TURNS, PEN sheets ---

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 185
$ws.Range("D2").Value = 10
$ws.Range("F2").Value = 23
$ws.Range("G2").Value = 41
$ws.Range("J2").Value = 30
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 190
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 146
$ws.Range("G3").Value = 59
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 60
$ws.Range("L3").Value = 380
$ws.Range("M3").Value = 256
$ws.Range("Q3").Value = 594

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 149
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 41
$ws.Range("J2").Value = 22
$ws.Range("N2").Value = 30
$ws.Range("O2").Value = 21
$ws.Range("C3").Value = 217
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 124
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 67
$ws.Range("L3").Value = 335
$ws.Range("M3").Value = 213
$ws.Range("Q3").Value = 575

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 89
$ws.Range("D2").Value = 54
$ws.Range("F2").Value = 457
$ws.Range("G2").Value = 439
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 185
$ws.Range("K2").Value = 175
$ws.Range("L2").Value = 118
$ws.Range("M2").Value = 89
$ws.Range("B3").Value = 72

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 16
$ws.Range("D2").Value = 16
$ws.Range("B3").Value = 12
$ws.Range("D3").Value = 6
$ws.Range("B4").Value = 2
